# Update cryptos list with latest prices/volumes from source feed
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.530.04"
$ws.Range("D3").Value = "3.255.71"
$ws.Range("E3").Value = "  +3.90%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'595.85"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").Value = "'141.62"
$ws.Range("E6").Value = "  +1.24%  "
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").Value = "3.250.70"
$ws.Range("E8").Value = "  +3.89%  "
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("E13").Value = "  -1.84%  "
$ws.Range("D14").Value = "'34.44"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").Value = "3.789.48"
$ws.Range("E15").Value = "  +3.80%  "
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "3.262.87"
$ws.Range("E17").Value = "  +4.11%  "
$ws.Range("D18").Value = "63.549.43"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "'6.79"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "'478.52"
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").Value = "'14.24"
$ws.Range("E21").Value = "  -1.92%  "
$ws.Range("D22").Value = "'0.732"
$ws.Range("E22").Value = "  +3.79%  "
$ws.Range("D23").Value = "'8.00"
$ws.Range("E23").Value = "  +4.13%  "
$ws.Range("D24").Value = "'83.77"
$ws.Range("E24").Value = "  -4.39%  "
$ws.Range("D25").Value = "'13.30"
$ws.Range("E25").Value = "  +1.43%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'2.75"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").Value = "'7.23"
$ws.Range("E28").Value = "  +4.20%  "
$ws.Range("D29").Value = "'8.11"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("E30").Value = "  +4.52%  "
$ws.Range("D31").Value = "'27.72"
$ws.Range("E31").Value = "  +2.03%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  -3.78%  "
$ws.Range("D34").Value = "'2.56"
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("E35").Value = "  -0.99%  "
$ws.Range("D36").Value = "'5.95"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("D37").Value = "'52.75"
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("D38").Value = "0.0₃0715"
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("D39").Value = "'0.0395"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "'423.48"
$ws.Range("E40").Value = "  -1.45%  "
$ws.Range("D41").Value = "3.003.33"
$ws.Range("E41").Value = "  +4.40%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.77"
$ws.Range("E42").Value = "  -2.65%  "
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").Value = "'8.39"
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("E44").Value = "  -7.47%  "
$ws.Range("D45").Value = "'0.266"
$ws.Range("E45").Value = "  +2.81%  "
$ws.Range("E46").Value = "  +1.32%  "
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").Value = "'2.35"
$ws.Range("E48").Value = "  -1.11%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'25.96"
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").Value = "'122.81"
$ws.Range("E51").Value = "  +1.96%  "
